$wb = $excel.ActiveWorkbook

# 1. Update the Date metadata value (Metadata sheet, row 8 = "Date" property)
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Cells.Item(8, 2).Value = "2025-07-21T08:56:18+00:00"

# 2. Remove the EXCL example row from the Concepts sheet.
#    Row 2 currently holds: Level=1, Code=EXCL, Display/Definition text.
#    Deleting it shifts the remaining ALT row up to become row 2.
$wsConcepts = $wb.Worksheets.Item("Concepts")
$wsConcepts.Rows.Item(2).Delete()
